$d = $word.ActiveDocument

$replacements = @(
    @{old = "314÷5=62, 4"; new = "283÷2=141, 1"},
    @{old = "217÷7=31, 0"; new = "576÷2=288, 0"},
    @{old = "969÷2=484, 1"; new = "178÷2=89, 0"},
    @{old = "766÷2=383, 0"; new = "447÷2=223, 1"},
    @{old = "875÷9=97, 2"; new = "215÷3=71, 2"},
    @{old = "723÷4=180, 3"; new = "909÷9=101, 0"},
    @{old = "975÷3=325, 0"; new = "191÷8=23, 7"},
    @{old = "322÷4=80, 2"; new = "483÷6=80, 3"},
    @{old = "468÷2=234, 0"; new = "442÷4=110, 2"},
    @{old = "744÷7=106, 2"; new = "423÷2=211, 1"},
    @{old = "497÷6=82, 5"; new = "549÷2=274, 1"},
    @{old = "657÷3=219, 0"; new = "104÷6=17, 2"},
    @{old = "935÷8=116, 7"; new = "194÷7=27, 5"},
    @{old = "587÷2=293, 1"; new = "270÷8=33, 6"},
    @{old = "783÷3=261, 0"; new = "956÷2=478, 0"},
    @{old = "435÷3=145, 0"; new = "593÷2=296, 1"},
    @{old = "115÷8=14, 3"; new = "216÷7=30, 6"},
    @{old = "903÷8=112, 7"; new = "262÷4=65, 2"},
    @{old = "607÷2=303, 1"; new = "935÷6=155, 5"},
    @{old = "111÷2=55, 1"; new = "539÷8=67, 3"},
    @{old = "160÷7=22, 6"; new = "374÷3=124, 2"},
    @{old = "129÷4=32, 1"; new = "293÷8=36, 5"},
    @{old = "416÷6=69, 2"; new = "811÷3=270, 1"},
    @{old = "878÷2=439, 0"; new = "183÷2=91, 1"},
    @{old = "651÷6=108, 3"; new = "753÷5=150, 3"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}

Write-Output "Done: applied $($replacements.Count) replacements"
